# RP2040 Board Ready For Review
# Fill in the two new resistor rows (10R fixed resistor + 0R jumper resistor)
# on the Resistors sheet, and update the active-sheet/selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resistors")

# --- Row 13: RES-FXD 10R 0603 chip resistor (fully sourced part) ---
$ws.Range("A13").Value = "RES-FXD"
$ws.Range("B13").Value = "10R"
$ws.Range("C13").Value = "1%"
$ws.Range("D13").Value = "0603"
$ws.Range("E13").Formula = "=_xlfn.XLOOKUP(D13,CaseCodes!A:A,CaseCodes!B:B)"
$ws.Range("F13").Value = "100ppm/°C"
$ws.Range("G13").Value = "100mW"
$ws.Range("H13").Value = "SMT"
$ws.Range("I13").Value = "Chip Resistor"
$ws.Range("J12:J13").Formula = "=UPPER(_xlfn.CONCAT(A12,""_"",B12,""_"",C12,""_"",G12,""_"",F12,""_"",D12,""("",E12,"")""))"
$ws.Range("K12:K13").Formula = "=_xlfn.CONCAT(B12,"" ±"",C12,"" "",G12,"" "",I12,"" "",F12,"" "",D12,""("",E12,"" Metric)"")"
$ws.Range("L13").Value = "EE490.SchLib"
$ws.Range("M13").Value = "EE490.PcbLib"
$ws.Range("N13").Value = "Resistor - Fixed - Generic"
$ws.Range("O13").Formula = "=_xlfn.XLOOKUP(D13,CaseCodes!A:A,CaseCodes!D:D)"
$ws.Range("P13").Value = "Bourns Inc."
$ws.Range("Q13").Value = "CR0603-FX-10R0ELF"
$ws.Range("R13").Value = "Active"
$ws.Range("S13").Value = "Bourns Inc."
$ws.Range("T13").Value = "CR0603AFX-10R0EAS"
$ws.Range("U13").Value = "Active"
$ws.Range("V13").Value = "Bourns Inc."
$ws.Range("W13").Value = "CMP0603AFX-10R0ELF"
$ws.Range("X13").Value = "Active"

# --- Row 12: RES-JMP 0R 0603 jumper resistor ---
$ws.Range("O12").Value = "RESISTOR_0603(1608)_NC_JUMPER"
$ws.Range("A12").Value = "RES-JMP"
$ws.Range("I12").Value = "Chip Resistor with Jumper"
$ws.Range("N12").Value = "Resistor - Fixed - Jumper"
$ws.Range("B12").Value = "0R"
$ws.Range("C12").Value = "1%"
$ws.Range("D12").Value = "0603"
$ws.Range("E12").Formula = "=_xlfn.XLOOKUP(D12,CaseCodes!A:A,CaseCodes!B:B)"
$ws.Range("F12").Value = "100ppm/°C"
$ws.Range("G12").Value = "100mW"
$ws.Range("H12").Value = "SMT"
$ws.Range("L12").Value = "EE490.SchLib"
$ws.Range("M12").Value = "EE490.PcbLib"

$excel.Calculate() | Out-Null

# --- View state: Resistors becomes the active/selected sheet, scrolled & with B14 selected ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("B14").Select() | Out-Null
